$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" column (C) serial date values for rows 2-7 move from
# 45184 (2023-09-15) to 45185 (2023-09-16).
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45185
}
